$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sigp")
$ws.Activate()

# Append 70 new rows (243:312) of id_department/id_criteria pairs
$data = New-Object 'object[,]' 70,2
$data[0,0] = 71; $data[0,1] = 237
$data[1,0] = 71; $data[1,1] = 238
$data[2,0] = 71; $data[2,1] = 239
$data[3,0] = 71; $data[3,1] = 240
$data[4,0] = 71; $data[4,1] = 241
$data[5,0] = 71; $data[5,1] = 242
$data[6,0] = 71; $data[6,1] = 243
$data[7,0] = 71; $data[7,1] = 244
$data[8,0] = 71; $data[8,1] = 245
$data[9,0] = 71; $data[9,1] = 246
$data[10,0] = 71; $data[10,1] = 247
$data[11,0] = 71; $data[11,1] = 248
$data[12,0] = 71; $data[12,1] = 249
$data[13,0] = 71; $data[13,1] = 250
$data[14,0] = 71; $data[14,1] = 251
$data[15,0] = 71; $data[15,1] = 252
$data[16,0] = 72; $data[16,1] = 253
$data[17,0] = 72; $data[17,1] = 254
$data[18,0] = 72; $data[18,1] = 255
$data[19,0] = 72; $data[19,1] = 256
$data[20,0] = 72; $data[20,1] = 257
$data[21,0] = 72; $data[21,1] = 258
$data[22,0] = 72; $data[22,1] = 259
$data[23,0] = 72; $data[23,1] = 260
$data[24,0] = 72; $data[24,1] = 261
$data[25,0] = 72; $data[25,1] = 262
$data[26,0] = 72; $data[26,1] = 263
$data[27,0] = 72; $data[27,1] = 264
$data[28,0] = 72; $data[28,1] = 265
$data[29,0] = 72; $data[29,1] = 266
$data[30,0] = 72; $data[30,1] = 267
$data[31,0] = 72; $data[31,1] = 268
$data[32,0] = 72; $data[32,1] = 269
$data[33,0] = 72; $data[33,1] = 270
$data[34,0] = 72; $data[34,1] = 271
$data[35,0] = 72; $data[35,1] = 272
$data[36,0] = 72; $data[36,1] = 273
$data[37,0] = 72; $data[37,1] = 274
$data[38,0] = 72; $data[38,1] = 275
$data[39,0] = 72; $data[39,1] = 276
$data[40,0] = 72; $data[40,1] = 277
$data[41,0] = 72; $data[41,1] = 278
$data[42,0] = 72; $data[42,1] = 279
$data[43,0] = 72; $data[43,1] = 280
$data[44,0] = 72; $data[44,1] = 281
$data[45,0] = 72; $data[45,1] = 282
$data[46,0] = 72; $data[46,1] = 283
$data[47,0] = 72; $data[47,1] = 284
$data[48,0] = 72; $data[48,1] = 285
$data[49,0] = 72; $data[49,1] = 286
$data[50,0] = 73; $data[50,1] = 287
$data[51,0] = 73; $data[51,1] = 288
$data[52,0] = 73; $data[52,1] = 289
$data[53,0] = 73; $data[53,1] = 290
$data[54,0] = 73; $data[54,1] = 291
$data[55,0] = 73; $data[55,1] = 292
$data[56,0] = 73; $data[56,1] = 293
$data[57,0] = 73; $data[57,1] = 294
$data[58,0] = 73; $data[58,1] = 295
$data[59,0] = 73; $data[59,1] = 296
$data[60,0] = 73; $data[60,1] = 297
$data[61,0] = 73; $data[61,1] = 298
$data[62,0] = 73; $data[62,1] = 299
$data[63,0] = 73; $data[63,1] = 300
$data[64,0] = 73; $data[64,1] = 301
$data[65,0] = 73; $data[65,1] = 302
$data[66,0] = 73; $data[66,1] = 303
$data[67,0] = 73; $data[67,1] = 304
$data[68,0] = 73; $data[68,1] = 305
$data[69,0] = 73; $data[69,1] = 306
$ws.Range("A243:B312").Value = $data

# Update selection to match the new view state
$ws.Range("B292").Select()

# Scroll viewport (best-effort; cosmetic window scroll position)
$excel.ActiveWindow.ScrollRow = 297
$excel.ActiveWindow.ScrollColumn = 1

# Adjust the saved workbook window width (cosmetic bookView size)
$excel.ActiveWindow.Width = 27945
